$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "estado de cuenta" detail rows (B16:G56) - part 1 of the new data,
# interleaving the two workers' periods (Isabel Cristina Gomez Martinez /
# Mailin Castro Banquez) in ascending period order.
$data = @(
    @(16, "CC", "1143396269", "ISABEL CRISTINA GOMEZ MARTINEZ", "2108", 36341, 908526),
    @(17, "CC", "1047451453", "MAILIN CASTRO BANQUEZ", "2109", 18170, 454263),
    @(18, "CC", "1143396269", "ISABEL CRISTINA GOMEZ MARTINEZ", "2109", 36341, 908526),
    @(19, "CC", "1047451453", "MAILIN CASTRO BANQUEZ", "2110", 18170, 454263),
    @(20, "CC", "1143396269", "ISABEL CRISTINA GOMEZ MARTINEZ", "2110", 36341, 908526),
    @(21, "CC", "1047451453", "MAILIN CASTRO BANQUEZ", "2111", 18170, 454263),
    @(22, "CC", "1143396269", "ISABEL CRISTINA GOMEZ MARTINEZ", "2111", 36341, 908526),
    @(23, "CC", "1047451453", "MAILIN CASTRO BANQUEZ", "2112", 18170, 454263),
    @(24, "CC", "1143396269", "ISABEL CRISTINA GOMEZ MARTINEZ", "2112", 36341, 908526),
    @(25, "CC", "1047451453", "MAILIN CASTRO BANQUEZ", "2201", 18170, 454263),
    @(26, "CC", "1143396269", "ISABEL CRISTINA GOMEZ MARTINEZ", "2201", 36341, 908526),
    @(27, "CC", "1047451453", "MAILIN CASTRO BANQUEZ", "2202", 18170, 454263),
    @(28, "CC", "1143396269", "ISABEL CRISTINA GOMEZ MARTINEZ", "2202", 36341, 908526),
    @(29, "CC", "1047451453", "MAILIN CASTRO BANQUEZ", "2203", 18170, 454263),
    @(30, "CC", "1143396269", "ISABEL CRISTINA GOMEZ MARTINEZ", "2203", 36341, 908526),
    @(31, "CC", "1047451453", "MAILIN CASTRO BANQUEZ", "2204", 18170, 454263),
    @(32, "CC", "1143396269", "ISABEL CRISTINA GOMEZ MARTINEZ", "2204", 36341, 908526),
    @(33, "CC", "1047451453", "MAILIN CASTRO BANQUEZ", "2205", 18170, 454263),
    @(34, "CC", "1143396269", "ISABEL CRISTINA GOMEZ MARTINEZ", "2205", 36341, 908526),
    @(35, "CC", "1047451453", "MAILIN CASTRO BANQUEZ", "2206", 18170, 454263),
    @(36, "CC", "1143396269", "ISABEL CRISTINA GOMEZ MARTINEZ", "2206", 36341, 908526),
    @(37, "CC", "1047451453", "MAILIN CASTRO BANQUEZ", "2207", 18170, 454263),
    @(38, "CC", "1143396269", "ISABEL CRISTINA GOMEZ MARTINEZ", "2207", 36341, 908526),
    @(39, "CC", "1047451453", "MAILIN CASTRO BANQUEZ", "2208", 18170, 454263),
    @(40, "CC", "1143396269", "ISABEL CRISTINA GOMEZ MARTINEZ", "2208", 36341, 908526),
    @(41, "CC", "1047451453", "MAILIN CASTRO BANQUEZ", "2209", 18170, 454263),
    @(42, "CC", "1143396269", "ISABEL CRISTINA GOMEZ MARTINEZ", "2209", 36341, 908526),
    @(43, "CC", "1047451453", "MAILIN CASTRO BANQUEZ", "2210", 18170, 454263),
    @(44, "CC", "1143396269", "ISABEL CRISTINA GOMEZ MARTINEZ", "2210", 36341, 908526),
    @(45, "CC", "1047451453", "MAILIN CASTRO BANQUEZ", "2211", 18170, 454263),
    @(46, "CC", "1143396269", "ISABEL CRISTINA GOMEZ MARTINEZ", "2211", 36341, 908526),
    @(47, "CC", "1047451453", "MAILIN CASTRO BANQUEZ", "2212", 18170, 454263),
    @(48, "CC", "1143396269", "ISABEL CRISTINA GOMEZ MARTINEZ", "2212", 36341, 908526),
    @(49, "CC", "1047451453", "MAILIN CASTRO BANQUEZ", "2301", 18170, 454263),
    @(50, "CC", "1143396269", "ISABEL CRISTINA GOMEZ MARTINEZ", "2301", 36341, 908526),
    @(51, "CC", "1047451453", "MAILIN CASTRO BANQUEZ", "2302", 18170, 454263),
    @(52, "CC", "1143396269", "ISABEL CRISTINA GOMEZ MARTINEZ", "2302", 36341, 908526),
    @(53, "CC", "1047451453", "MAILIN CASTRO BANQUEZ", "2303", 18170, 454263),
    @(54, "CC", "1143396269", "ISABEL CRISTINA GOMEZ MARTINEZ", "2303", 36341, 908526),
    @(55, "CC", "1047451453", "MAILIN CASTRO BANQUEZ", "2304", 16353, 454263),
    @(56, "CC", "1143396269", "ISABEL CRISTINA GOMEZ MARTINEZ", "2304", 32707, 908526)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
